$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the team's Win/Loss/Tie record (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered style used by the other header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record (68 wins, 94 losses, 0 ties) for every player row
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 94   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
